# Atualizado por script em 11-11-2023 20:45
#
# 1) Swap the match data (columns F:V) between rows 84 and 85 (the
#    Mladost vs Radnicki Nis and Vojvodina vs Radnik matches had their
#    row order reversed). Columns A:E (Indice, pais, torneio, temporada,
#    data_partida) are unaffected.
# 2) Append three new match rows (105-107) at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap F84:V84 with F85:V85 --------------------------------------
$swapCols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

$row84 = @{}
$row85 = @{}
foreach ($c in $swapCols) {
    $row84[$c] = $ws.Range($c + "84").Value2
    $row85[$c] = $ws.Range($c + "85").Value2
}
foreach ($c in $swapCols) {
    $ws.Range($c + "84").Value = $row85[$c]
    $ws.Range($c + "85").Value = $row84[$c]
}

# --- 2) Append new rows 105, 106, 107 -----------------------------------
# Copy formatting (styles) from the last existing row (104) down onto the
# three new rows first, then overwrite the values.
$ws.Range("A104:V104").Copy($ws.Range("A105:V105"))
$ws.Range("A104:V104").Copy($ws.Range("A106:V106"))
$ws.Range("A104:V104").Copy($ws.Range("A107:V107"))

# Row 105
$ws.Range("A105").Value = 104
$ws.Range("B105").Value = "serbia"
$ws.Range("C105").Value = "super-liga"
$ws.Range("D105").Value = "2023-2024"
$ws.Range("E105").Value = 45241.66666666666
$ws.Range("F105").Value = "Radnik"
$ws.Range("G105").Value = 1
$ws.Range("H105").Value = "Napredak"
$ws.Range("I105").Value = 2
$ws.Range("J105").Value = 2.39
$ws.Range("K105").Value = "09/11/2023 09:13"
$ws.Range("L105").Value = 2.31
$ws.Range("M105").Value = "11/11/2023 15:59"
$ws.Range("N105").Value = 2.87
$ws.Range("O105").Value = "09/11/2023 09:13"
$ws.Range("P105").Value = 2.76
$ws.Range("Q105").Value = "11/11/2023 15:59"
$ws.Range("R105").Value = 2.97
$ws.Range("S105").Value = "09/11/2023 09:13"
$ws.Range("T105").Value = 3.61
$ws.Range("U105").Value = "11/11/2023 15:59"
$ws.Range("V105").Value = "https://www.betexplorer.com/football/serbia/super-liga/radnik-surdulica-napredak/MDe2IpJa/"

# Row 106
$ws.Range("A106").Value = 105
$ws.Range("B106").Value = "serbia"
$ws.Range("C106").Value = "super-liga"
$ws.Range("D106").Value = "2023-2024"
$ws.Range("E106").Value = 45241.75
$ws.Range("F106").Value = "Zeleznicar Pancevo"
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = "Crvena zvezda"
$ws.Range("I106").Value = 2
$ws.Range("J106").Value = 6.73
$ws.Range("K106").Value = "09/11/2023 09:13"
$ws.Range("L106").Value = 21.51
$ws.Range("M106").Value = "11/11/2023 17:58"
$ws.Range("N106").Value = 5.01
$ws.Range("O106").Value = "09/11/2023 09:13"
$ws.Range("P106").Value = 8.359999999999999
$ws.Range("Q106").Value = "11/11/2023 17:58"
$ws.Range("R106").Value = 1.3
$ws.Range("S106").Value = "09/11/2023 09:13"
$ws.Range("T106").Value = 1.11
$ws.Range("U106").Value = "11/11/2023 17:08"
$ws.Range("V106").Value = "https://www.betexplorer.com/football/serbia/super-liga/zeleznicar-pancevo-crvena-zvezda/lzSbH4Y5/"

# Row 107
$ws.Range("A107").Value = 106
$ws.Range("B107").Value = "serbia"
$ws.Range("C107").Value = "super-liga"
$ws.Range("D107").Value = "2023-2024"
$ws.Range("E107").Value = 45241.77083333334
$ws.Range("F107").Value = "IMT Novi Beograd"
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = "Mladost"
$ws.Range("I107").Value = 1
$ws.Range("J107").Value = 1.83
$ws.Range("K107").Value = "09/11/2023 09:13"
$ws.Range("L107").Value = 1.76
$ws.Range("M107").Value = "11/11/2023 18:29"
$ws.Range("N107").Value = 3.33
$ws.Range("O107").Value = "09/11/2023 09:13"
$ws.Range("P107").Value = 3.7
$ws.Range("Q107").Value = "11/11/2023 18:29"
$ws.Range("R107").Value = 3.71
$ws.Range("S107").Value = "09/11/2023 09:13"
$ws.Range("T107").Value = 4.3
$ws.Range("U107").Value = "11/11/2023 18:29"
$ws.Range("V107").Value = "https://www.betexplorer.com/football/serbia/super-liga/imt-novi-beograd-mladost-lucani/rNdbJQ3g/"
